# Weekly update: insert two new price records after the existing row 17
# (before row 18), pushing the former rows 18-41 down to rows 20-43.
#
# New row 18 is a variant of the old row 18 (Asterix / Provincia de Melipilla)
# with an updated date, quality and prices.
# New row 19 is a variant of the old row 19 (Rodeo / Región de Los Lagos)
# but the variety/quality/prices/origin now match a "Cardinal" /
# "Provincia de Melipilla" record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 18 and 19; Excel shifts rows 18:41 down to 20:43
# and copies the formatting (incl. the date style on column D) from the
# row immediately above, just like a normal Excel "Insert Rows" would.
$ws.Rows("18:19").Insert()

# ---- Row 18 (carries forward A/B/C/H/O from the old row 18) ----
$ws.Range("A18").Value = 1
$ws.Range("B18").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C18").Value = "Arica y Parinacota"
$ws.Range("D18").Value = "2021-08-19"
$ws.Range("E18").Value = 15
$ws.Range("F18").Value = 100114001
$ws.Range("G18").Value = "Papa"
$ws.Range("H18").Value = "Asterix"
$ws.Range("I18").Value = "1a (cosecha lavada)"
$ws.Range("J18").Value = 1000
$ws.Range("K18").Value = 9000
$ws.Range("L18").Value = 10000
$ws.Range("M18").Value = 9500
$ws.Range("N18").Value = "$/malla 25 kilos"
$ws.Range("O18").Value = "Provincia de Melipilla"
$ws.Range("P18").Value = 380
$ws.Range("Q18").Value = 25
$ws.Range("R18").Value = "Hortaliza"

# ---- Row 19 (carries forward A/B/C from the old row 19) ----
$ws.Range("A19").Value = 1
$ws.Range("B19").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C19").Value = "Arica y Parinacota"
$ws.Range("D19").Value = "2021-08-19"
$ws.Range("E19").Value = 15
$ws.Range("F19").Value = 100114001
$ws.Range("G19").Value = "Papa"
$ws.Range("H19").Value = "Cardinal"
$ws.Range("I19").Value = "1a (cosecha)"
$ws.Range("J19").Value = 1000
$ws.Range("K19").Value = 11000
$ws.Range("L19").Value = 12000
$ws.Range("M19").Value = 11500
$ws.Range("N19").Value = "$/saco 25 kilos"
$ws.Range("O19").Value = "Provincia de Melipilla"
$ws.Range("P19").Value = 460
$ws.Range("Q19").Value = 25
$ws.Range("R19").Value = "Hortaliza"
